$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22, pushing existing rows 22..118 down to 23..119.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new weekly record.
$ws.Range("A22").Value = 8
$ws.Range("B22").Value = "Terminal La Palmera de La Serena"
$ws.Range("C22").Value = "Coquimbo"
$ws.Range("D22").Value = 44972
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = 100112030
$ws.Range("G22").Value = "Poroto granado"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 37000
$ws.Range("L22").Value = 38000
$ws.Range("M22").Value = 37500
$ws.Range("N22").Value = "$/malla 25 kilos"
$ws.Range("O22").Value = "Provincia del Elquí"
$ws.Range("P22").Value = 1500
$ws.Range("Q22").Value = 25
$ws.Range("R22").Value = "Hortaliza"
